# Weekly price-list update: a new weekly record is inserted above the
# existing row 73 ("Feria Lagunitas de Puerto Montt" - Espárragos), shifting
# the subsequent rows (old 73-85) down to 74-86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 73 (existing rows 73.. shift down by one).
$ws.Rows(73).Insert()

# Populate the newly inserted row 73 with the new weekly record.
$ws.Range("A73").Value = 4
$ws.Range("B73").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C73").Value = "Los Lagos"
$ws.Range("D73").Value = 45244
$ws.Range("E73").Value = 10
$ws.Range("F73").Value = 300000000
$ws.Range("G73").Value = "Espárragos"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 400
$ws.Range("K73").Value = 2000
$ws.Range("L73").Value = 2000
$ws.Range("M73").Value = 2000
$ws.Range("N73").Value = "$/kilo"
$ws.Range("O73").Value = "Provincia de Linares"
$ws.Range("P73").Value = 2000
$ws.Range("Q73").Value = 1
$ws.Range("R73").Value = "Hortaliza"
